$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.059661865234375
$ws.Range("B1").Value = 2.259135007858276
$ws.Range("C1").Value = 2.354539632797241
$ws.Range("D1").Value = 3.045063018798828
$ws.Range("E1").Value = 2.967832565307617
